# ---------------------------------------------------------------------------
# Rebuild Sheet1 of the electric-bridge workbook: add the missing S/R2/R1-R2/
# Rx columns, add R0' uncertainty columns D and I, populate the new K:P
# "uncertainty" block, border/format the data table, resize columns, and
# restore the saved selection + page setup.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row text (shared strings). New unique strings "S", "R2(欧)",
#    "R1/R2", "Rx(欧)" get appended to sharedStrings.xml the first time they
#    are used below.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "R1(欧)"
$ws.Range("B1").Value = "R0(欧)"
$ws.Range("C1").Value = "R0'(欧)"
$ws.Range("D1").Value = "S"
$ws.Range("F1").Value = "E(V)"
$ws.Range("G1").Value = "R0(欧)"
$ws.Range("H1").Value = "R0'(欧)"
$ws.Range("I1").Value = "S"
$ws.Range("K1").Value = "E(V)"
$ws.Range("L1").Value = "R1(欧)"
$ws.Range("M1").Value = "R2(欧)"
$ws.Range("N1").Value = "R1/R2"
$ws.Range("O1").Value = "R0(欧)"
$ws.Range("P1").Value = "Rx(欧)"

# ---------------------------------------------------------------------------
# 2) Establish the cellXfs styles in the exact order the workbook ends up
#    needing them (2=border+center, 3=border, 4=numFmt .0, 5=numFmt .0;red,
#    6=border+numFmt.0+center, 7=border+numFmt.0, 8=border+numFmt.0;red+
#    center, 9=border+numFmt.0;red) by touching one representative cell per
#    style first, then broadcasting the same formatting to the remaining
#    cells that need it (which reuses the style already created).
# ---------------------------------------------------------------------------
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("D2").NumberFormat = "0.0_ "
$ws.Range("L2").NumberFormat = "0.0_);[Red]\(0.0\)"
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").NumberFormat = "0.0_ "
$ws.Range("B2").Borders.LineStyle = 1
$ws.Range("B2").NumberFormat = "0.0_ "
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").Borders.LineStyle = 1
$ws.Range("L1").NumberFormat = "0.0_);[Red]\(0.0\)"
$ws.Range("L2").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3) Broadcast borders to every header/data cell of the table.
# ---------------------------------------------------------------------------
# Header row (border + center align) -> style 2
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("I1").Borders.LineStyle = 1
$ws.Range("K1").Borders.LineStyle = 1
$ws.Range("N1").Borders.LineStyle = 1
$ws.Range("P1").Borders.LineStyle = 1

# Header row (border + numFmt .0 + center align) -> style 6
$ws.Range("C1").Borders.LineStyle = 1
$ws.Range("C1").NumberFormat = "0.0_ "
$ws.Range("G1").Borders.LineStyle = 1
$ws.Range("G1").NumberFormat = "0.0_ "
$ws.Range("H1").Borders.LineStyle = 1
$ws.Range("H1").NumberFormat = "0.0_ "
$ws.Range("M1").Borders.LineStyle = 1
$ws.Range("M1").NumberFormat = "0.0_ "
$ws.Range("O1").Borders.LineStyle = 1
$ws.Range("O1").NumberFormat = "0.0_ "

# Data rows (border only) -> style 3
$ws.Range("A2:A6").Borders.LineStyle = 1
$ws.Range("D2:D6").Borders.LineStyle = 1
$ws.Range("F2:F5").Borders.LineStyle = 1
$ws.Range("I2:I5").Borders.LineStyle = 1
$ws.Range("K2:K4").Borders.LineStyle = 1
$ws.Range("N2:N4").Borders.LineStyle = 1
$ws.Range("P2:P4").Borders.LineStyle = 1

# Data rows (border + numFmt .0) -> style 7
$ws.Range("B2:B6").Borders.LineStyle = 1
$ws.Range("B2:B6").NumberFormat = "0.0_ "
$ws.Range("C2:C6").Borders.LineStyle = 1
$ws.Range("C2:C6").NumberFormat = "0.0_ "
$ws.Range("G2:G5").Borders.LineStyle = 1
$ws.Range("G2:G5").NumberFormat = "0.0_ "
$ws.Range("H2:H5").Borders.LineStyle = 1
$ws.Range("H2:H5").NumberFormat = "0.0_ "
$ws.Range("M2:M4").Borders.LineStyle = 1
$ws.Range("M2:M4").NumberFormat = "0.0_ "
$ws.Range("O2:O4").Borders.LineStyle = 1
$ws.Range("O2:O4").NumberFormat = "0.0_ "

# Data rows (border + numFmt .0;red) -> style 9
$ws.Range("L2:L4").Borders.LineStyle = 1
$ws.Range("L2:L4").NumberFormat = "0.0_);[Red]\(0.0\)"

# ---------------------------------------------------------------------------
# 4) Cell values / formulas.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 900
$ws.Range("B2").Value = 61469
$ws.Range("C2").Value = 61669
$ws.Range("D2").Value = 122.938
$ws.Range("F2").Value = 2.5
$ws.Range("G2").Value = 61583
$ws.Range("H2").Value = 61883
$ws.Range("I2").Value = 82.110667000000007
$ws.Range("K2").Value = 2.5
$ws.Range("L2").Value = 900
$ws.Range("M2").Value = 9000
$ws.Range("N2").Formula = "=L2/M2"
$ws.Range("O2").Value = 61199
$ws.Range("P2").Formula = "=O2/M2*L2"

$ws.Range("A3").Value = 1000
$ws.Range("B3").Value = 61719
$ws.Range("C3").Value = 61919
$ws.Range("D3").Value = 123.438
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 61383
$ws.Range("H3").Value = 61983
$ws.Range("I3").Value = 40.921999999999997
$ws.Range("K3").Value = 2.5
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 30000
$ws.Range("N3").Formula = "=L3/M3"
$ws.Range("O3").Value = 625.4
$ws.Range("P3").Formula = "=O3/M3*L3"

$ws.Range("A4").Value = 2000
$ws.Range("B4").Value = 61509
$ws.Range("C4").Value = 61809
$ws.Range("D4").Value = 82.012
$ws.Range("F4").Value = 1.5
$ws.Range("G4").Value = 61483
$ws.Range("H4").Value = 61933
$ws.Range("I4").Value = 54.651555999999999
$ws.Range("K4").Value = 2.5
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 20000
$ws.Range("N4").Formula = "=L4/M4"
$ws.Range("O4").Value = 3274
$ws.Range("P4").Formula = "=O4/M4*L4"

$ws.Range("A5").Value = 3000
$ws.Range("B5").Value = 61609
$ws.Range("C5").Value = 61905
$ws.Range("D5").Value = 83.255405409999995
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 61199
$ws.Range("H5").Value = 62199
$ws.Range("I5").Value = 24.479600000000001

$ws.Range("A6").Value = 5000
$ws.Range("B6").Value = 61543
$ws.Range("C6").Value = 61883
$ws.Range("D6").Value = 72.403529410000004

# ---------------------------------------------------------------------------
# 5) Column widths (best effort - the COM layer quantizes ColumnWidth to
#    1/7-character pixel steps, so these are the closest achievable values).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 6.142857142857143
$ws.Columns.Item(2).ColumnWidth = 7.0
$ws.Columns.Item(3).ColumnWidth = 7.142857142857143
$ws.Columns.Item(4).ColumnWidth = 6.714285714285714
$ws.Columns.Item(5).ColumnWidth = 3.2857142857142856
$ws.Columns.Item(6).ColumnWidth = 5.142857142857143
$ws.Columns.Item(7).ColumnWidth = 7.428571428571429
$ws.Columns.Item(8).ColumnWidth = 7.285714285714286
$ws.Columns.Item(9).ColumnWidth = 6.142857142857143
$ws.Columns.Item(10).ColumnWidth = 3.2857142857142856
$ws.Columns.Item(11).ColumnWidth = 4.428571428571429
$ws.Columns.Item(12).ColumnWidth = 7.428571428571429
$ws.Columns.Item(13).ColumnWidth = 7.428571428571429
$ws.Columns.Item(14).ColumnWidth = 5.571428571428571
$ws.Columns.Item(15).ColumnWidth = 8.285714285714286
$ws.Columns.Item(16).ColumnWidth = 6.857142857142857

# ---------------------------------------------------------------------------
# 6) Page setup (paper size / orientation) and selected cell.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("H14").Select()
